$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update hours logged on 2019-02-13 (row 25) from 5.5 to 6.5
$ws.Range("B25").Value = 6.5

# Append additional work notes to the task description in D25
$ws.Range("D25").Value = "Indie Project: MVP ERD, populated lookup tables, tried to figure mysqldump (saved creation & insert files separately for now), created User entity and draft dao, copied in SessionFactoryProvider"

# Update the current selection to D26 (matches author's cursor position on save)
$ws.Range("D26").Select()
